# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Estado de Cuenta" table (rows 16-27) previously listed each worker's
# six overdue periods together (DEIVIS 1905..1812, then MARLIN 1905..1812).
# This update reorganizes the table so the periods run in ascending order
# (1812, 1901, 1902, 1903, 1904, 1905) with both workers reported for each
# period, and adds MARLIN's new period rows as part 1 of the new account
# statement data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: C = N Doc Trabajador, D = Nombre Trabajador, E = Periodo Mora, F = Valor Mora

$rows = @(
    @{ Row = 16; Doc = "73185071";   Nombre = "DEIVIS RAFAEL BARON SALCEDO";      Periodo = "1812"; Valor = 31249 },
    @{ Row = 17; Doc = "1047458640"; Nombre = "MARLIN DEL CARMEN MARRUGO VELEZ";  Periodo = "1812"; Valor = 31249 },
    @{ Row = 18; Doc = "73185071";   Nombre = "DEIVIS RAFAEL BARON SALCEDO";      Periodo = "1901"; Valor = 31249 },
    @{ Row = 19; Doc = "1047458640"; Nombre = "MARLIN DEL CARMEN MARRUGO VELEZ";  Periodo = "1901"; Valor = 31249 },
    @{ Row = 20; Doc = "73185071";   Nombre = "DEIVIS RAFAEL BARON SALCEDO";      Periodo = "1902"; Valor = 31249 },
    @{ Row = 21; Doc = "1047458640"; Nombre = "MARLIN DEL CARMEN MARRUGO VELEZ";  Periodo = "1902"; Valor = 31249 },
    @{ Row = 22; Doc = "73185071";   Nombre = "DEIVIS RAFAEL BARON SALCEDO";      Periodo = "1903"; Valor = 31249 },
    @{ Row = 23; Doc = "1047458640"; Nombre = "MARLIN DEL CARMEN MARRUGO VELEZ";  Periodo = "1903"; Valor = 31249 },
    @{ Row = 24; Doc = "73185071";   Nombre = "DEIVIS RAFAEL BARON SALCEDO";      Periodo = "1904"; Valor = 31249 },
    @{ Row = 25; Doc = "1047458640"; Nombre = "MARLIN DEL CARMEN MARRUGO VELEZ";  Periodo = "1904"; Valor = 31249 },
    @{ Row = 26; Doc = "73185071";   Nombre = "DEIVIS RAFAEL BARON SALCEDO";      Periodo = "1905"; Valor = 20833 },
    @{ Row = 27; Doc = "1047458640"; Nombre = "MARLIN DEL CARMEN MARRUGO VELEZ";  Periodo = "1905"; Valor = 20833 }
)

foreach ($r in $rows) {
    $ws.Range("C$($r.Row)").Value = $r.Doc
    $ws.Range("D$($r.Row)").Value = $r.Nombre
    $ws.Range("E$($r.Row)").Value = $r.Periodo
    $ws.Range("F$($r.Row)").Value = $r.Valor
}
